# Updated symbol list on Sat Feb 11 07:48:45 UTC 2023 with GitHub Actions
# All cells in columns D (Price) and E (Volume(1h)) hold numeric-looking text,
# so force the number format to Text ("@") before assigning the string value;
# otherwise Excel would auto-convert them to real numbers/percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '308.13'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.31%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '40.67'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '2.14%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.128'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.16%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07617'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-1.26%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.617'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-0.73%'
$ws.Range("B7").Value = 'MXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9024'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '2.63%'
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.429'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-0.35%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1103'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '10.60%'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '0.96%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09086'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '1.90%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04179'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-5.10%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.1050'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.50%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001257'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.13%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005799'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.54%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.353'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.00%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.254'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.33%'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.93%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.579'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-6.39%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1364'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '1.90%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-5.85%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04067'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-2.15%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001222'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.85%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004088'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.82%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001301'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '6.75%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02378'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '1.14%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05172'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '0.40%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007769'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-2.12%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.006772'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '6.87%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.001952'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '0.94%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008763'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '3.17%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3327'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '8.42%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00007014'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '8.03%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.16%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.03098'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '377.55%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.004201'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-40.03%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002102'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.16%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002002'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.16%'
